# Auto-generated edit script: apply scheduled price-refresh update
# to Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 54770.3
$ws.Range("I69").Value = 3241.1428
$ws.Range("K69").Value = 9723.428400000001
$ws.Range("M69").Value = -8849.428400000001

$ws.Range("H72").Value = 54770.3
$ws.Range("I72").Value = 3241.1428
$ws.Range("K72").Value = 29170.2852
$ws.Range("M72").Value = -24802.2852

$ws.Range("H86").Value = 3857.3333
$ws.Range("I86").Value = 3520.6
$ws.Range("K86").Value = 3520.6
$ws.Range("M86").Value = -2397.6

$ws.Range("H87").Value = 54652
$ws.Range("J87").Value = 54652
$ws.Range("L87").Value = 54652
$ws.Range("N87").Value = -57148

$ws.Range("H89").Value = 3857.3333
$ws.Range("I89").Value = 3520.6
$ws.Range("K89").Value = 17603
$ws.Range("M89").Value = -11987

$ws.Range("H90").Value = 54652
$ws.Range("J90").Value = 54652
$ws.Range("L90").Value = 163956
$ws.Range("N90").Value = -176436

$ws.Range("H111").Value = 1039
$ws.Range("I111").Value = 1039
$ws.Range("K111").Value = 3117
$ws.Range("M111").Value = -50

$ws.Range("H138").Value = 5067.607
$ws.Range("I138").Value = 4850
$ws.Range("J138").Value = 5084.346
$ws.Range("K138").Value = 14550
$ws.Range("L138").Value = 15253.038
$ws.Range("M138").Value = -9410
$ws.Range("N138").Value = -25533.038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16102.223
$ws.Range("I32").Value = 14808.077
$ws.Range("K32").Value = 14808.077
$ws.Range("M32").Value = -14521.077

$ws.Range("H88").Value = 1583.3334
$ws.Range("I88").Value = 1066.6666
$ws.Range("J88").Value = 2100
$ws.Range("K88").Value = 1066.6666
$ws.Range("L88").Value = 2100
$ws.Range("M88").Value = -660.6666
$ws.Range("N88").Value = -2912

$ws.Range("H91").Value = 1583.3334
$ws.Range("I91").Value = 1066.6666
$ws.Range("J91").Value = 2100
$ws.Range("K91").Value = 1066.6666
$ws.Range("L91").Value = 2100
$ws.Range("M91").Value = 337.3334
$ws.Range("N91").Value = -4908

$ws.Range("H110").Value = 2758.8948
$ws.Range("I110").Value = 1367.1333
$ws.Range("K110").Value = 1367.1333
$ws.Range("M110").Value = 677.8667

$ws.Range("H122").Value = 3562
$ws.Range("I122").Value = 3202.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9607.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -7157.5
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2786
$ws.Range("I94").Value = 3296.5
$ws.Range("J94").Value = 999.25
$ws.Range("K94").Value = 3296.5
$ws.Range("L94").Value = 999.25
$ws.Range("M94").Value = -2845.5
$ws.Range("N94").Value = -1901.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2333.5293
$ws.Range("I31").Value = 1841.8
$ws.Range("J31").Value = 3036
$ws.Range("K31").Value = 1841.8
$ws.Range("L31").Value = 3036
$ws.Range("M31").Value = -1546.8
$ws.Range("N31").Value = -3626

$ws.Range("H34").Value = 2333.5293
$ws.Range("I34").Value = 1841.8
$ws.Range("J34").Value = 3036
$ws.Range("K34").Value = 1841.8
$ws.Range("L34").Value = 3036
$ws.Range("M34").Value = -1639.8
$ws.Range("N34").Value = -3440

$ws.Range("H60").Value = 35572
$ws.Range("J60").Value = 47909
$ws.Range("L60").Value = 47909
$ws.Range("N60").Value = -48931

$ws.Range("H107").Value = 2227
$ws.Range("I107").Value = 1518
$ws.Range("K107").Value = 1518
$ws.Range("M107").Value = 402

$ws.Range("H118").Value = 61499.5
$ws.Range("J118").Value = 61499.5
$ws.Range("L118").Value = 61499.5
$ws.Range("N118").Value = -64813.5

$ws.Range("H131").Value = 35800.5
$ws.Range("J131").Value = 35800.5
$ws.Range("L131").Value = 35800.5
$ws.Range("N131").Value = -45880.5

$ws.Range("H141").Value = 424489.25
$ws.Range("J141").Value = 424489.25
$ws.Range("L141").Value = 424489.25
$ws.Range("N141").Value = -434849.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 263177
$ws.Range("I4").Value = 333334.2
$ws.Range("J4").Value = 87.5
$ws.Range("K4").Value = 1000002.6
$ws.Range("L4").Value = 262.5
$ws.Range("M4").Value = -999890.6000000001
$ws.Range("N4").Value = -486.5

$ws.Range("H12").Value = 329.6
$ws.Range("I12").Value = 5.5
$ws.Range("K12").Value = 16.5
$ws.Range("M12").Value = 156.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H11").Value = 10000000
$ws.Range("I11").Value = 10000000
$ws.Range("K11").Value = 10000000
$ws.Range("M11").Value = -9999861

$ws.Range("H132").Value = 5496.4443
$ws.Range("I132").Value = 4094
$ws.Range("K132").Value = 12282
$ws.Range("M132").Value = -9752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 21233.334
$ws.Range("I24").Value = 15480
$ws.Range("J24").Value = 50000
$ws.Range("K24").Value = 15480
$ws.Range("L24").Value = 50000
$ws.Range("M24").Value = -15137
$ws.Range("N24").Value = -50686

$ws.Range("H46").Value = 1299.875
$ws.Range("I46").Value = 799.75
$ws.Range("J46").Value = 1800
$ws.Range("K46").Value = 799.75
$ws.Range("L46").Value = 1800
$ws.Range("M46").Value = -611.75
$ws.Range("N46").Value = -2176

$ws.Range("H82").Value = 2819.6667
$ws.Range("I82").Value = 2443.6
$ws.Range("J82").Value = 4700
$ws.Range("K82").Value = 2443.6
$ws.Range("L82").Value = 4700
$ws.Range("M82").Value = -2082.6
$ws.Range("N82").Value = -5422

$ws.Range("H85").Value = 2819.6667
$ws.Range("I85").Value = 2443.6
$ws.Range("J85").Value = 4700
$ws.Range("K85").Value = 2443.6
$ws.Range("L85").Value = 4700
$ws.Range("M85").Value = -1195.6
$ws.Range("N85").Value = -7196

$ws.Range("H93").Value = 1983.1666
$ws.Range("I93").Value = 1720.091
$ws.Range("K93").Value = 1720.091
$ws.Range("M93").Value = -472.0909999999999

$ws.Range("H131").Value = 39325.25
$ws.Range("J131").Value = 39325.25
$ws.Range("L131").Value = 39325.25
$ws.Range("N131").Value = -49405.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1683333.4
$ws.Range("I3").Value = 1683333.4
$ws.Range("K3").Value = 1683333.4
$ws.Range("M3").Value = -1683219.4

$ws.Range("H14").Value = 3134.6667
$ws.Range("I14").Value = 400
$ws.Range("J14").Value = 4502
$ws.Range("K14").Value = 400
$ws.Range("L14").Value = 4502
$ws.Range("M14").Value = -232
$ws.Range("N14").Value = -4838

$ws.Range("H62").Value = 7880.375
$ws.Range("I62").Value = 4400
$ws.Range("J62").Value = 8377.571
$ws.Range("K62").Value = 4400
$ws.Range("L62").Value = 8377.571
$ws.Range("M62").Value = -3776
$ws.Range("N62").Value = -9625.571

$ws.Range("H65").Value = 7880.375
$ws.Range("I65").Value = 4400
$ws.Range("J65").Value = 8377.571
$ws.Range("K65").Value = 22000
$ws.Range("L65").Value = 41887.855
$ws.Range("M65").Value = -18880
$ws.Range("N65").Value = -48127.855

$ws.Range("H107").Value = 695.6429000000001
$ws.Range("I107").Value = 595.1429000000001
$ws.Range("J107").Value = 796.1429000000001
$ws.Range("K107").Value = 1785.4287
$ws.Range("L107").Value = 2388.4287
$ws.Range("M107").Value = 134.5712999999998
$ws.Range("N107").Value = -6228.4287

$ws.Range("H126").Value = 1651.6522
$ws.Range("I126").Value = 1473.5
$ws.Range("J126").Value = 1689.1578
$ws.Range("K126").Value = 4420.5
$ws.Range("L126").Value = 5067.4734
$ws.Range("M126").Value = -1950.5
$ws.Range("N126").Value = -10007.4734
